# Data edit: the "psf HDB Avg" column (F) for four districts had been left
# at 0 and is now populated with the computed average price-per-square-foot
# values (rows 7, 10, 25 and 27 on Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F7").Value  = 464.57922710399998
$ws.Range("F10").Value = 421.135501836
$ws.Range("F25").Value = 324.82616766899997
$ws.Range("F27").Value = 358.15189464000002

# Reflect the view state the workbook was saved with: Sheet1 active/selected,
# zoomed to 100% with the last selection sitting on F27 (the cell most
# recently edited).
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F27").Select()
